{"js": "// Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific resume bullet\n// paragraphs, matching the author's \"Implement quantitative metrics\n// highlighting\" commit.\n//\n// Strategy: identify each target paragraph by a unique, stable substring of\n// its original text, then \u2014 scoped to just that paragraph \u2014 search() for\n// each metric substring and flip its run(s) to bold + the accent color.\n// Word/Office.js automatically splits the run(s) at the match boundaries,\n// which reproduces the diff's run-splitting exactly.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: a substring unique to the paragraph (used to locate it) and\n// the ordered list of metric substrings within that paragraph to bold+color.\nconst TARGETS = [\n  {\n    match:\n      \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    match:\n      \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"],\n  },\n  {\n    match: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    match:\n      \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    match: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    match: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of TARGETS) {\n  const para = paragraphs.items.find((p) => p.text === target.match);\n  if (!para) continue;\n\n  for (const metric of target.metrics) {\n    const found = para.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    for (const run of found.items) {\n      run.font.bold = true;\n      run.font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific resume bullet\n# paragraphs, matching the author's \"Implement quantitative metrics\n# highlighting\" commit.\n#\n# Strategy: walk $d.Paragraphs, identify each target paragraph by a unique,\n# stable substring of its original text, then \u2014 scoped to just that\n# paragraph's [Start,End) span \u2014 run Find.Execute for each metric substring\n# and bold + color the found range. Word automatically splits the run(s) at\n# the match boundaries, reproducing the diff's run-splitting exactly.\n\n$d = $word.ActiveDocument\n\n# BGR-ordered OLE color value for #2C3E50 (Word/VBA Font.Color is BGR, not RGB).\n$HighlightColor = 0x50 * 65536 + 0x3E * 256 + 0x2C\n\n# Each entry: a substring unique to the paragraph (used to locate it) and the\n# ordered list of metric substrings within that paragraph to bold + color.\n$Targets = @(\n    @{\n        Match   = \"Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Match   = \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Metrics = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Match   = \"Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Match   = \"Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Metrics = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Match   = \"Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Match   = \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\nforeach ($target in $Targets) {\n    foreach ($p in $d.Paragraphs) {\n        $text = $p.Range.Text.TrimEnd()\n        # Strip the leading bullet glyph + space (if present) before comparing.\n        if ($text.StartsWith([string]([char]0x2022) + \" \")) {\n            $text = $text.Substring(2)\n        }\n\n        if ($text -eq $target.Match) {\n            $pStart = $p.Range.Start\n            $pEnd = $p.Range.End\n\n            foreach ($needle in $target.Metrics) {\n                $r = $d.Range($pStart, $pEnd)\n                $found = $r.Find.Execute($needle, $true)\n                if ($found) {\n                    $r.Font.Bold = 1\n                    $r.Font.Color = $HighlightColor\n                }\n            }\n\n            break\n        }\n    }\n}\n"}
